$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update build number in C2
$ws.Range("C2").Value = 7821

# Add new changelog entry row, copying the style of A4
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5").Value = "Changed lightmap 1.21.6+ (Vanilla)"

# Move selection to A6
$ws.Range("A6").Select()
